$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emirates")

# Update flight search dates (Search Flight_001 row)
# Leading apostrophe preserves the existing "stored as text" cell
# formatting (quote-prefix style) instead of re-evaluating it.
$ws.Range("D2").Value = "'11X2023"
$ws.Range("E2").Value = "'20X2023"

# Update flight search dates (Search Flight_002 row)
$ws.Range("D3").Value = "'10X2023"
$ws.Range("E3").Value = "'12X2023"
